$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (30 "characters" - ColumnWidth 29.2 rounds to a stored width of 30)
for ($c = 1; $c -le 9; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 29.2
}

# Header row
$ws.Cells.Item(1, 1).Value = "Broker Name"
$ws.Cells.Item(1, 2).Value = "Amount In Hand"
$ws.Cells.Item(1, 3).Value = "No. of Customers"
$ws.Cells.Item(1, 4).Value = "Settled Customers"
$ws.Cells.Item(1, 5).Value = "Total Customers"
$ws.Cells.Item(1, 6).Value = "Loan Given"
$ws.Cells.Item(1, 7).Value = "Loan Settled"
$ws.Cells.Item(1, 8).Value = "Loan till date"
$ws.Cells.Item(1, 9).Value = "Broker Address"

# Row 3 (divesh)
$ws.Cells.Item(3, 1).Value = "divesh"
$ws.Cells.Item(3, 2).Value = 9979578
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 20422
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 20422
$ws.Cells.Item(3, 9).Value = "muthu street"

# Row 4 (nilesh) - Amount In Hand stored as text (leading apostrophe forces text storage)
$ws.Cells.Item(4, 1).Value = "nilesh"
$ws.Cells.Item(4, 2).Value = "'120000"
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = "muthu street"

$ws.Range("A1").Select()
